$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Range("D19").Value = "아기 있는 집에는 에몬스홈 그란데 가죽소파 4인"
$ws.Range("E19").Value = "https://kkokkilkon.tistory.com/183"

# Row 20
$ws.Range("D20").Value = "[python 독학] 9. 상대경로, 절대경로 - 파일 read, write할 때 필수"
$ws.Range("E20").Value = "https://ai-creator.tistory.com/539"

# Row 23
$ws.Range("D23").Value = "파이참pycharm에서 2개 이상의 multiple projects 돌리는 방법"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2745"

# Row 51
$ws.Range("D51").Value = "[세이버메트릭스] 타율, 출루율, OPS, RC 중 무엇이 가장 득점 생산과 연관 있을까?"
$ws.Range("E51").Value = "https://bskyvision.com/1132"
